$wb = $excel.ActiveWorkbook

# The "Croatia" test-data sheet is a near duplicate of "Turkey" - duplicate the
# Turkey worksheet (placing the copy immediately after it) and then localize
# the two market-specific cells plus the tab name.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Copy($null, $turkey) | Out-Null

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2415/T2416/T2417"

# Restore Turkey's own selection (select-all) and leave it as a non-active
# tab, then make Croatia the active tab with its own selection - matching
# the recorded end-user state after adding the new sheet.
$turkey.Activate()
$turkey.Cells.Select() | Out-Null

$croatia.Activate()
$croatia.Range("C18").Select() | Out-Null

Write-Host "Added Croatia worksheet after Turkey"
